$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for columns B:E in the data rows so that numeric-looking
# strings (e.g. "1.0000", "27.159.88") are preserved exactly as text and not
# auto-converted to numbers/dates by Excel.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.159.88'
$ws.Range("E2").Value = '  -3.15%  '

$ws.Range("D3").Value = '1.716.55'
$ws.Range("E3").Value = '  -3.43%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = '311.54'
$ws.Range("E5").Value = '  -5.44%  '

$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("D7").Value = '0.4808'
$ws.Range("E7").Value = '  +7.00%  '

$ws.Range("D8").Value = '0.3459'
$ws.Range("E8").Value = '  -3.00%  '

$ws.Range("D9").Value = '42.77'
$ws.Range("E9").Value = '  +1.41%  '

$ws.Range("D10").Value = '0.07298'
$ws.Range("E10").Value = '  -2.11%  '

$ws.Range("D11").Value = '1.048'
$ws.Range("E11").Value = '  -5.42%  '

$ws.Range("E12").Value = '  -0.16%  '

$ws.Range("D13").Value = '19.95'
$ws.Range("E13").Value = '  -5.04%  '

$ws.Range("D14").Value = '5.880'
$ws.Range("E14").Value = '  -3.00%  '

$ws.Range("D15").Value = '1.715.67'
$ws.Range("E15").Value = '  -3.55%  '

$ws.Range("D16").Value = '6.883'
$ws.Range("E16").Value = '  -5.35%  '

$ws.Range("D17").Value = '89.24'
$ws.Range("E17").Value = '  -4.78%  '

$ws.Range("D18").Value = '0.00001042'
$ws.Range("E18").Value = '  -2.06%  '

$ws.Range("D19").Value = '0.06368'
$ws.Range("E19").Value = '  -1.04%  '

$ws.Range("D20").Value = '0.9999'
$ws.Range("E20").Value = '  -0.05%  '

$ws.Range("D21").Value = '16.54'
$ws.Range("E21").Value = '  -3.63%  '

$ws.Range("D22").Value = '5.654'
$ws.Range("E22").Value = '  -2.60%  '

$ws.Range("D23").Value = '27.216.43'
$ws.Range("E23").Value = '  -3.02%  '

$ws.Range("D24").Value = '10.85'
$ws.Range("E24").Value = '  -4.42%  '

$ws.Range("D25").Value = '2.088'
$ws.Range("E25").Value = '  -1.72%  '

$ws.Range("D26").Value = '152.21'
$ws.Range("E26").Value = '  -5.93%  '

$ws.Range("D27").Value = '19.73'
$ws.Range("E27").Value = '  -3.31%  '

$ws.Range("D28").Value = '1.909.21'
$ws.Range("E28").Value = '  -3.74%  '

$ws.Range("D29").Value = '2.105'
$ws.Range("E29").Value = '  -2.70%  '

$ws.Range("D30").Value = '120.46'
$ws.Range("E30").Value = '  -3.67%  '

$ws.Range("D31").Value = '1.024'
$ws.Range("E31").Value = '  -7.94%  '

$ws.Range("D32").Value = '0.09263'
$ws.Range("E32").Value = '  +0.57%  '

$ws.Range("D33").Value = '3.581'

$ws.Range("D34").Value = '5.359'
$ws.Range("E34").Value = '  -6.07%  '

$ws.Range("D35").Value = '0.02205'
$ws.Range("E35").Value = '  -3.89%  '

$ws.Range("D36").Value = '0.05936'
$ws.Range("E36").Value = '  -4.55%  '

$ws.Range("D37").Value = '11.13'
$ws.Range("E37").Value = '  -6.44%  '

$ws.Range("E38").Value = '  -5.12%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '4.776'
$ws.Range("E39").Value = '  -4.52%  '

$ws.Range("B40").Value = 'WEMIXTOKEN'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").Value = '1.422'
$ws.Range("E40").Value = '  +1.70%  '

$ws.Range("D41").Value = '0.9994'
$ws.Range("E41").Value = '  -0.14%  '

$ws.Range("D42").Value = '0.5957'
$ws.Range("E42").Value = '  -6.00%  '

$ws.Range("D43").Value = '1.100'
$ws.Range("E43").Value = '  -7.30%  '

$ws.Range("D44").Value = '7.492'
$ws.Range("E44").Value = '  -5.44%  '

$ws.Range("D45").Value = '12.73'
$ws.Range("E45").Value = '  -4.64%  '

$ws.Range("D46").Value = '3.591'
$ws.Range("E46").Value = '  -4.43%  '

$ws.Range("D47").Value = '0.5636'
$ws.Range("E47").Value = '  -4.85%  '

$ws.Range("D48").Value = '118.82'
$ws.Range("E48").Value = '  -3.45%  '

$ws.Range("D49").Value = '1.848'
$ws.Range("E49").Value = '  -5.89%  '

$ws.Range("D50").Value = '0.06650'
$ws.Range("E50").Value = '  -3.64%  '

$ws.Range("D51").Value = '1.087'
$ws.Range("E51").Value = '  -5.42%  '

Write-Host "cryptos list updated"
